# Append two newly-scraped Lancers listings to the top portion of the
# "ランサーズ" table, shifting the existing rows down (matching the site's
# "most recent first" ordering) and re-stamping every row's "取得日時"
# (fetched-at) timestamp to the new run's time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-05 01:53:48"

# --- 1) Make room: insert the two brand-new rows -------------------------
# New row #1 lands at sheet row 2 (top of the data, newest first).
$ws.Rows.Item(2).Insert()
# New row #2 lands at sheet row 7 (post-shift position from the diff).
$ws.Rows.Item(7).Insert()

# --- 2) Fill the two newly-inserted rows ----------------------------------
$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = "SaaSビジネスにおける「バーティカル(垂直型)」展開の横スライド可能なAIシステムの開発です"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5485911"
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("G2").Value = 383
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

$ws.Range("A7").Value = $newTimestamp
$ws.Range("B7").Value = "【急募】クリックポスト自動発行ツール開発依頼"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5485895"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("G7").Value = 123
$ws.Range("H7").Value = "◆ツール,開発"

# --- 3) Re-stamp the "取得日時" column for every other (shifted) row ------
# (rows 3-6 and 8-20 hold the pre-existing listings, now shifted down by
# the two inserts above; only column A changes for them.)
foreach ($r in 3..6) {
    $ws.Range("A" + $r).Value = $newTimestamp
}
foreach ($r in 8..20) {
    $ws.Range("A" + $r).Value = $newTimestamp
}

# --- 4) Rebuild the hyperlinks collection ---------------------------------
# Row-insert in this host does not slide hyperlink anchors along with the
# cells, so drop the stale collection and re-add one hyperlink per F-cell,
# top to bottom, to match the final F2:F20 layout.
$ws.Range("A1").Hyperlinks.Delete()

$urls = @{
    2  = "https://www.lancers.jp/work/detail/5485911"
    3  = "https://www.lancers.jp/work/detail/5455098"
    4  = "https://www.lancers.jp/work/detail/5445159"
    5  = "https://www.lancers.jp/work/detail/5445154"
    6  = "https://www.lancers.jp/work/detail/5485147"
    7  = "https://www.lancers.jp/work/detail/5485895"
    8  = "https://www.lancers.jp/work/detail/5485816"
    9  = "https://www.lancers.jp/work/detail/5485476"
    10 = "https://www.lancers.jp/work/detail/5485054"
    11 = "https://www.lancers.jp/work/detail/5485362"
    12 = "https://www.lancers.jp/work/detail/5484588"
    13 = "https://www.lancers.jp/work/detail/5481888"
    14 = "https://www.lancers.jp/work/detail/5485460"
    15 = "https://www.lancers.jp/work/detail/5485304"
    16 = "https://www.lancers.jp/work/detail/5485174"
    17 = "https://www.lancers.jp/work/detail/5485814"
    18 = "https://www.lancers.jp/work/detail/5485734"
    19 = "https://www.lancers.jp/work/detail/5485785"
    20 = "https://www.lancers.jp/work/detail/5485506"
}

foreach ($r in 2..20) {
    $ws.Hyperlinks.Add($ws.Range("F" + $r), $urls[$r])
    $ws.Range("F" + $r).Style = "Hyperlink"
}

"ok"
